# Section 11 - Configure Shopping Cart API with new methods
#
# Insert a new paragraph, right after the existing
# "add-migration AddCouponTodb" paragraph, containing:
#   "add-migration" <run> + " " <run> + "seedCouponDatabase" <run,
#   wrapped in spellStart/spellEnd proofErr marks>
# matching the same run formatting (Cascadia Mono, size 19, black)
# already used by the sibling command paragraphs in this section.

$d = $word.ActiveDocument

# Locate the paragraph whose text is "add-migration AddCouponTodb" --
# this is the anchor the new paragraph must be inserted after.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*AddCouponTodb*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'AddCouponTodb' paragraph to anchor the insert."
}

# Create a new, empty paragraph immediately after it (inherits the same
# paragraph mark formatting as the anchor paragraph).
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()

# Build the OOXML for the new paragraph's content: three runs
# ("add-migration", " ", "seedCouponDatabase") each with identical
# run formatting, with spellStart/spellEnd proofErr marks bracketing
# the last run -- mirrors the structure already used elsewhere in the
# document for these inline command paragraphs.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
  '<w:pPr>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/>' +
      '<w:color w:val="000000"/>' +
      '<w:sz w:val="19"/>' +
      '<w:szCs w:val="19"/>' +
    '</w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/>' +
      '<w:color w:val="000000"/>' +
      '<w:sz w:val="19"/>' +
      '<w:szCs w:val="19"/>' +
    '</w:rPr>' +
    '<w:t>add-migration</w:t>' +
  '</w:r>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/>' +
      '<w:color w:val="000000"/>' +
      '<w:sz w:val="19"/>' +
      '<w:szCs w:val="19"/>' +
    '</w:rPr>' +
    '<w:t xml:space="preserve"> </w:t>' +
  '</w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/>' +
      '<w:color w:val="000000"/>' +
      '<w:sz w:val="19"/>' +
      '<w:szCs w:val="19"/>' +
    '</w:rPr>' +
    '<w:t>seedCouponDatabase</w:t>' +
  '</w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

# Replace the whole (empty) new paragraph's range with our fully-formed
# paragraph XML, so no stray empty run is left behind.
$newPara.Range.InsertXML($xml)

Write-Output "Inserted 'add-migration seedCouponDatabase' paragraph after 'AddCouponTodb'."
